$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.011.94'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.550.85'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.95%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '287.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3915'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.61%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3198'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.26'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07268'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.093'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('E13').Value = '  -7.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.617'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.88%  '
$ws.Range('E15').Value = '  -2.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001124'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.548.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.87%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06580'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '83.46'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.95%  '
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.276'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.70'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.21'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.023.38'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.353'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.27%  '
$ws.Range('E26').Value = '  -3.82%  '
$ws.Range('E27').Value = '  -2.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.54'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.57%  '
$ws.Range('E29').Value = '  -1.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.721.84'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '118.69'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.052'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.667'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08331'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.162'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.599'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -16.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06146'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02259'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.101'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.214'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2062'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.45%  '
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.57'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5795'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.21'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.711'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5547'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.92%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '117.89'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.01%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.892'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.136'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.60%  '
$ws.Range('E51').Value = '  -4.30%  '
